# Updated symbol list on Thu Dec 29 19:05:32 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.20"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "19"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "24.18"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "19"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.284"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "19"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05777"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "19"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.504"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "19"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.143"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "19"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8121"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "19"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8620"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "19"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1377"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "19"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06959"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "19"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03152"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "19"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02919"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "19"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09384"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "19"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.763"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "19"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001525"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "19"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04671"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "19"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005979"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "19"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006165"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "19"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "19"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.004649"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "19"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00006099"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "19"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.506"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "19"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.150"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "19"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "19"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "19"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "19"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0002331"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "19"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "19"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "19"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "19"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "19"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "19"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "19"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "19"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "19"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "19"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "19"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "19"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03709"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "19"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006279"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "19"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1057"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "19"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003400"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "19"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007724"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "19"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "19"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "19"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4399"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "19"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002410"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "19"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "19"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "19"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "19"
